$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.8779
$ws.Range("E6").Value = 12.59600000000001
$ws.Range("E7").Value = 12.11439999999999
$ws.Range("D8").Value = -8.860599999999987
$ws.Range("E8").Value = 12.40319999999999
$ws.Range("A12").Value = -22.86360000000002
$ws.Range("D12").Value = -8.272700000000004
$ws.Range("D14").Value = -8.664000000000001
$ws.Range("E19").Value = 12.857
$ws.Range("E21").Value = 12.48579999999999
$ws.Range("D22").Value = -7.932999999999993
$ws.Range("E24").Value = 12.97379999999999
